# Update "想去人数" (interested-count) figures on the 展览 and 全部类型 sheets
# to match the newly generated gh-pages output (commit 456a3b4).

$wb = $excel.ActiveWorkbook

# --- 展览 (sheet1): rows keyed by F-column old value -> new value ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Cells.Item(3, 6).Value  = 73      # 张家港·上元节AF 汉服花灯动漫展: 72 -> 73
$ws1.Cells.Item(9, 6).Value  = 269     # 苏州·世纪幻想动漫游戏展2.0: 268 -> 269
$ws1.Cells.Item(11, 6).Value = 10118   # 【会员购严选】苏州·Anime LIVE 国际动漫品牌博览会: 10105 -> 10118
$ws1.Cells.Item(13, 6).Value = 263     # 苏州·绘时国乙1.0-秩序之外: 262 -> 263
$ws1.Cells.Item(15, 6).Value = 630     # 苏州·梦幻岛 国乙主题文化展（日夜场） 梦幻岛之约3.0: 628 -> 630
$ws1.Cells.Item(16, 6).Value = 11775   # 昆山·第十二届理想乡动漫游戏展: 11772 -> 11775
$ws1.Cells.Item(17, 6).Value = 12151   # 苏州·第十七届 I COME ACG 动漫品牌博览会: 12143 -> 12151

# --- 全部类型 (sheet4): same events, offset by one row further down ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Cells.Item(3, 6).Value  = 73      # 张家港·上元节AF 汉服花灯动漫展: 72 -> 73
$ws4.Cells.Item(10, 6).Value = 269     # 苏州·世纪幻想动漫游戏展2.0: 268 -> 269
$ws4.Cells.Item(12, 6).Value = 10118   # 【会员购严选】苏州·Anime LIVE 国际动漫品牌博览会: 10105 -> 10118
$ws4.Cells.Item(14, 6).Value = 263     # 苏州·绘时国乙1.0-秩序之外: 262 -> 263
$ws4.Cells.Item(16, 6).Value = 630     # 苏州·梦幻岛 国乙主题文化展（日夜场） 梦幻岛之约3.0: 628 -> 630
$ws4.Cells.Item(17, 6).Value = 11775   # 昆山·第十二届理想乡动漫游戏展: 11772 -> 11775
$ws4.Cells.Item(18, 6).Value = 12151   # 苏州·第十七届 I COME ACG 动漫品牌博览会: 12143 -> 12151
